# ModelRuns_RTP2025.xlsx — add three new 2023 test runs (IPA_51/52/53)
# that lower Work_Transit_Hesitance while testing different rail/ferry
# hesitance values, inserted right after the existing 2023 rows (before
# the 2025/2035/2050 blocks).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 70:72 - formatting is inherited from the row
# above (row 69), and every row/formula below shifts down automatically.
$ws.Rows("70:72").Insert()

# Column B ("directory") - filled top to bottom first.
$ws.Range("B70").Value = "2023_TM160_IPA_51"
$ws.Range("B71").Value = "2023_TM160_IPA_52"
$ws.Range("B72").Value = "2023_TM160_IPA_53"

# Column F ("description") - filled top to bottom next.
$ws.Range("F70").Value = "AOC=16.21, with rail/ferry hes =70"
$ws.Range("F71").Value = "AOC=16.21, with rail/ferry hes =65"
$ws.Range("F72").Value = "AOC=16.21, with rail/ferry hes =60"

# Asana links - filled bottom row (72) up to top row (70), matching the
# order the shared-string table was populated by the original author.
$ws.Range("L72").Value = "https://app.asana.com/0/1204085012544660/1206685870182016/f"
$ws.Range("L71").Value = "https://app.asana.com/0/1204085012544660/1206685870182013/f"
$ws.Range("L70").Value = "https://app.asana.com/0/1204085012544660/1206699871295080/f"

# "model_machine" (column K) - still to-be-determined for all three runs.
$ws.Range("K70").Value = "tbd"
$ws.Range("K71").Value = "tbd"
$ws.Range("K72").Value = "tbd"

# --- Remaining columns, same for all three new rows -----------------------
$ws.Range("A70:A72").Value = 2023
$ws.Range("C70:C72").Value = "RTP2025_IP"
$ws.Range("D70:D72").Value = "Base year"
$ws.Range("G70:G72").Value = "petrale"
$ws.Range("H70:H72").Value = "n/a"
$ws.Range("I70:I72").Value = "current"
$ws.Range("J70:J72").Value = "BlueprintNetworks_v13\net_2023_Blueprint"
$ws.Range("M70").Value = 16.21
$ws.Range("M71").Value = 16.21
$ws.Range("M72").Value = 16.21
$ws.Range("N70:N72").Value = "na"
$ws.Range("O70:O72").Value = "na"
$ws.Range("P70:P72").Value = 1.04
$ws.Range("Q70:Q72").Value = 0.94
$ws.Range("R70").Value = 90
$ws.Range("R71").Value = 80
$ws.Range("R72").Value = 70
$ws.Range("S70:S72").Value = 0
$ws.Range("T70:T72").Value = 75

# Match the selection left by the author in the saved file.
$ws.Range("K72").Select()
